$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-11 (data trimmed to only the 3 most recent entries)
$ws.Rows("5:11").Delete()

# Update row 2: Salary, 500, 2025-05-23
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 45800.12527777778

# Update row 3: Gift, 100, 2025-05-12
$ws.Range("A3").Value = "Gift"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 45789.12527777778

# Update row 4: Freelance, 250, 2025-05-11
$ws.Range("A4").Value = "Freelance"
$ws.Range("B4").Value = 250
$ws.Range("C4").Value = 45788.12527777778
